# Apply the "#5: insurance, claim, debt, investment done" edit.
#
# This extends the "債務" (Debt, sheet index 5) and "事業投資" (Business
# investment, sheet index 6) worksheets with the full set of metadata
# columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) that the other sheets in this
# workbook already carry, and fixes up the header rows which previously
# held stray data values instead of proper column labels.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $addr, $val) {
    # New header cells need to look like the existing bold / bordered /
    # centered header style (cellXfs index 1) used by B1:G1.
    $r = $ws.Range($addr)
    $r.Value = $val
    $r.Font.Bold = $true
    $r.Borders.LineStyle = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
}

function Set-DataCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-DataCellText($ws, $addr, $val) {
    # Plain "$ws.Range($addr).Value = $val" lets Excel auto-detect values
    # that look like dates (e.g. "2012-04-12") and silently convert them
    # into date serial numbers with a new number-format style. Routing
    # the literal text through a TEXT() formula, then pasting back just
    # the computed value, keeps it as a genuine string cell without
    # touching any cell formatting / styles.
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '=TEXT("' + $val + '","@")'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------
# Sheet 5: 債務 (Debt)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Row 1 (header) - B1/C1/D1/E1/F1/G1 previously held leftover data values;
# replace them with the correct column headers and extend through N1.
Set-HeaderCell $ws5 "B1" "species"
Set-HeaderCell $ws5 "C1" "debtor"
Set-HeaderCell $ws5 "D1" "owner"
Set-HeaderCell $ws5 "E1" "total"
Set-HeaderCell $ws5 "F1" "register_date"
Set-HeaderCell $ws5 "G1" "register_reason"
Set-HeaderCell $ws5 "H1" "property_category"
Set-HeaderCell $ws5 "I1" "category"
Set-HeaderCell $ws5 "J1" "date"
Set-HeaderCell $ws5 "K1" "legislator_name"
Set-HeaderCell $ws5 "L1" "legislator_id"
Set-HeaderCell $ws5 "M1" "source_file"
Set-HeaderCell $ws5 "N1" "index"

# Row 2 (data) - B2:G2 already hold the correct values; only add H2:N2.
Set-DataCell $ws5 "H2" "debt"
Set-DataCell $ws5 "I2" "normal"
Set-DataCellText $ws5 "J2" "2012-04-12"
Set-DataCell $ws5 "K2" "盧嘉辰"
Set-DataCell $ws5 "L2" 1715
Set-DataCell $ws5 "M2" "tmp79201"
Set-DataCell $ws5 "N2" 95

# ---------------------------------------------------------------------
# Sheet 6: 事業投資 (Business investment)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Row 1 (header) - B1/C1/D1/E1/F1/G1 previously held leftover data values;
# replace them with the correct column headers and extend through N1.
Set-HeaderCell $ws6 "B1" "owner"
Set-HeaderCell $ws6 "C1" "company"
Set-HeaderCell $ws6 "D1" "address"
Set-HeaderCell $ws6 "E1" "total"
Set-HeaderCell $ws6 "F1" "register_date"
Set-HeaderCell $ws6 "G1" "register_reason"
Set-HeaderCell $ws6 "H1" "property_category"
Set-HeaderCell $ws6 "I1" "category"
Set-HeaderCell $ws6 "J1" "date"
Set-HeaderCell $ws6 "K1" "legislator_name"
Set-HeaderCell $ws6 "L1" "legislator_id"
Set-HeaderCell $ws6 "M1" "source_file"
Set-HeaderCell $ws6 "N1" "index"

# Row 2 (data) - B2:G2 already hold the correct values; only add H2:N2.
Set-DataCell $ws6 "H2" "investment"
Set-DataCell $ws6 "I2" "normal"
Set-DataCellText $ws6 "J2" "2012-04-12"
Set-DataCell $ws6 "K2" "盧嘉辰"
Set-DataCell $ws6 "L2" 1715
Set-DataCell $ws6 "M2" "tmp79201"
Set-DataCell $ws6 "N2" 100

# Row 3 (data) - B3:G3 already hold the correct values; only add H3:N3.
Set-DataCell $ws6 "H3" "investment"
Set-DataCell $ws6 "I3" "normal"
Set-DataCellText $ws6 "J3" "2012-04-12"
Set-DataCell $ws6 "K3" "盧嘉辰"
Set-DataCell $ws6 "L3" 1715
Set-DataCell $ws6 "M3" "tmp79201"
Set-DataCell $ws6 "N3" 101
